# Remove canceled contestants from backup data for accurate record keeping.
#
# 1) The "Canceled Assignments" sheet (and its one data row, tied to the
#    contestant whose assignment was canceled) is removed outright.
# 2) On the Contestants sheet the canceled contestant (Peter Adamidis, the
#    current row 2) drops to the bottom of the table, the other two
#    contestants shift up one row, and everyone's AttendingWith list is
#    rewired to match the new arrangement.

$wb = $excel.ActiveWorkbook

# --- 1. Drop the "Canceled Assignments" sheet ------------------------------
$canceled = $wb.Worksheets.Item("Canceled Assignments")
$canceled.Delete()

# --- 2. Re-home the Contestants rows ---------------------------------------
$ws = $wb.Worksheets.Item("Contestants")

# The "Location" (G) and "Rating" (H) columns are blank for some
# contestants, but whether the cell is present at all (vs. entirely
# missing from the row) varies, and that presence/absence needs to move
# along with each contestant. Settle that *before* touching any other
# columns, while the sheet still holds its original values:
#   - G4 (currently "Melbourne") must become present-but-blank, like G2
#     currently is. Clear it first so the destination starts absent, same
#     as the H3 case below (copying onto an already non-blank cell is a
#     no-op in this engine when the source is blank).
#   - H3 (currently absent) must become present-but-blank, like H4
#     currently is.
#   - H2 (currently present-but-blank) must become absent.
$ws.Cells.Item(4, 7).Value() = ""
$ws.Range("G2").Copy($ws.Range("G4"))
$ws.Range("H4").Copy($ws.Range("H3"))
$ws.Cells.Item(2, 8).Value() = ""

# Row 2 <- Kathleen Reynolds (previously row 3)
# Column F (Phone) is "498086080" for every contestant, so it is left
# untouched below rather than re-written (keeps its original text form).
# Column G is handled above.
$ws.Cells.Item(2, 1).Value()  = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$ws.Cells.Item(2, 2).Value()  = "Kathleen Reynolds"
$ws.Cells.Item(2, 3).Value()  = 33
$ws.Cells.Item(2, 4).Value()  = "Not Specified"
$ws.Cells.Item(2, 5).Value()  = "kathleenmonicareynolds@gmail.com"
$ws.Cells.Item(2, 7).Value()  = "Footscray"
$ws.Cells.Item(2, 9).Value()  = "available"
$ws.Cells.Item(2, 10).Value() = "Peter Adamidis, Felicity Parker-Hill"
$ws.Cells.Item(2, 11).Value() = "5fe641da-4067-49a7-bae7-e63413b3e404"
$ws.Cells.Item(2, 12).Value() = "N"
$ws.Cells.Item(2, 13).Value() = "N/A"

# Row 3 <- Felicity Parker-Hill (previously row 4)
$ws.Cells.Item(3, 1).Value()  = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$ws.Cells.Item(3, 2).Value()  = "Felicity Parker-Hill"
$ws.Cells.Item(3, 3).Value()  = 27
$ws.Cells.Item(3, 4).Value()  = "Not Specified"
$ws.Cells.Item(3, 5).Value()  = "felicity.parkerhill@endemolshine.com.au"
$ws.Cells.Item(3, 7).Value()  = "Melbourne"
$ws.Cells.Item(3, 9).Value()  = "available"
$ws.Cells.Item(3, 10).Value() = "Peter Adamidis, Kathleen Reynolds"
$ws.Cells.Item(3, 11).Value() = "5fe641da-4067-49a7-bae7-e63413b3e404"
$ws.Cells.Item(3, 12).Value() = "N"
$ws.Cells.Item(3, 13).Value() = "N/A"

# Row 4 <- Peter Adamidis (previously row 2, the canceled contestant)
$ws.Cells.Item(4, 1).Value()  = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$ws.Cells.Item(4, 2).Value()  = "Peter Adamidis"
$ws.Cells.Item(4, 3).Value()  = 34
$ws.Cells.Item(4, 4).Value()  = "Not Specified"
$ws.Cells.Item(4, 5).Value()  = "peter.adamidis@gmail.com"
$ws.Cells.Item(4, 9).Value()  = "available"
$ws.Cells.Item(4, 10).Value() = "Kathleen Reynolds, Felicity Parker-Hill"
$ws.Cells.Item(4, 11).Value() = "5fe641da-4067-49a7-bae7-e63413b3e404"
$ws.Cells.Item(4, 12).Value() = "Y"
$ws.Cells.Item(4, 13).Value() = "Broken Leg"
